$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "dSF" (column F) values per repulled data / mean calculation
$updates = @{
    2  = -6
    5  = -2
    6  = -5
    8  = 3
    10 = -3
    11 = -7
    12 = -3
    13 = -6
    15 = -1
    17 = 6
    18 = -6
    20 = -4
    22 = -1
    24 = 5
    28 = 2
    29 = 8
    30 = 1
    31 = -5
    33 = 1
    34 = 1
    35 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
